$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): blank out A1, then drop the bold/bordered/centered style ---
$ws.Range("A1").Value = ""
$ws.Range("A1:P1").ClearFormats()

# --- Corrected data cleaning for pre/post/total fixation data ---

# Row 3 - Revisit count
$ws.Range("C3").Value = 25
$ws.Range("I3").Value = 2
$ws.Range("L3").Value = 18
$ws.Range("M3").Value = 23
$ws.Range("P3").Value = 1

# Row 4 - Fixation count
$ws.Range("C4").Value = 104
$ws.Range("I4").Value = 4
$ws.Range("L4").Value = 47
$ws.Range("M4").Value = 95
$ws.Range("P4").Value = 2

# Row 5 - Dwell time (ms)
$ws.Range("C5").Value = 22272.74
$ws.Range("I5").Value = 1251.37
$ws.Range("L5").Value = 10995.88
$ws.Range("M5").Value = 34673.09
$ws.Range("P5").Value = 834.28

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 0.4
$ws.Range("C6").Value = 25.47
$ws.Range("E6").Value = 1.58
$ws.Range("F6").Value = 0.31
$ws.Range("H6").Value = 0.4
$ws.Range("I6").Value = 1.43
$ws.Range("J6").Value = 3.64
$ws.Range("L6").Value = 12.57
$ws.Range("M6").Value = 39.65
$ws.Range("N6").Value = 0.17
$ws.Range("P6").Value = 0.95

# Row 7 - Fixation duration (ms)
$ws.Range("C7").Value = 214.16
$ws.Range("I7").Value = 312.84
$ws.Range("L7").Value = 233.95
$ws.Range("M7").Value = 364.98
$ws.Range("P7").Value = 417.14

# --- Row 10 was an erroneous extra blank row; remove it entirely ---
$ws.Rows.Item(10).Delete()
